$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row heights for row 1 and row 2 ---
$ws.Rows.Item(1).RowHeight = 41.25
$ws.Rows.Item(2).RowHeight = 15

# --- New column R (year 2021) values, keyed by row number ---
$values = @{
  4  = 2021
  6  = 88.796593100633856
  7  = 86.908583391486388
  8  = 89.680106631122953
  9  = 95.775910364145659
  10 = 96.517042279754136
  11 = 90.311530128242666
  12 = 90.746324915190343
  13 = 90.894107952204379
  14 = 81.065680730752504
  15 = 85.088888888888889
  17 = 93.37839883628321
  18 = 93.091416608513612
  19 = 94.815061646117954
  20 = 100.53781512605042
  21 = 100.33525796237662
  22 = 93.78989283832054
  23 = 95.401432340746325
  24 = 92.308748798242007
  25 = 89.338842975206617
  26 = 87.955555555555549
  28 = 89.631204460036727
  29 = 89.204466154919743
  30 = 84.751749416861045
  31 = 96.201680672268907
  32 = 95.567144719687093
  33 = 91.330444457457389
  34 = 91.368262344515642
  35 = 92.345373803964662
  36 = 88.660287081339717
  37 = 84.944444444444443
}

# Column R mirrors column Q's formatting all the way down (rows 3-38),
# same as every other year column added to this table.
for ($r = 3; $r -le 38; $r++) {
  $null = $ws.Range("Q$r").Copy()
  $null = $ws.Range("R$r").PasteSpecial(-4122)

  if ($values.ContainsKey($r)) {
    $ws.Range("R$r").Value = $values[$r]
  }
}

# --- Update the stored selection to R3 ---
$null = $ws.Range("R3").Select()
